$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("E1").Value = "Nombre Completo"
$ws.Range("J1").Value = "Filial Operacion"

# Sort data rows (A2:L13) by column H (Fecha de Operacion) descending
$rng = $ws.Range("A1:L13")
$rng.Sort($ws.Range("H1"), 2, $null, $null, 1, $null, 1, 1)
